# Apply transformer parameter additions to the prediction parameters sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for transformer params
$ws.Range("K1").Value = "d_model"
$ws.Range("L1").Value = "num_layers"

# Update existing row-2 values
$ws.Range("C2").Value = 23
$ws.Range("D2").Value = 110
$ws.Range("F2").Value = 0.0005
$ws.Range("H2").Value = 10

# New row-2 values for the transformer params
$ws.Range("K2").Value = 16
$ws.Range("L2").Value = 2

# Update selection to match target state
$ws.Range("F5").Select()
